$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Table 1 ("Patient Characteristics") - Gender section
#   Row 2 = "Gender" (header row, p-value cell)
#   Row 3 = "    Female"
#   Row 4 = "    Male"
#   Row 5 = "    Unknown/Unspecified"  -> entire row removed
# ------------------------------------------------------------------
$t1 = $d.Tables.Item(1)

# Gender p-value
$t1.Cell(2, 5).Range.Text = "0.2"

# Female row
$t1.Cell(3, 2).Range.Text = "30 (28%)"
$t1.Cell(3, 3).Range.Text = "18 (35%)"
$t1.Cell(3, 4).Range.Text = "12 (22%)"

# Male row
$t1.Cell(4, 2).Range.Text = "76 (72%)"
$t1.Cell(4, 3).Range.Text = "34 (65%)"
$t1.Cell(4, 4).Range.Text = "42 (78%)"

# Remove the "Unknown/Unspecified" row entirely
$t1.Rows.Item(5).Delete()

# ------------------------------------------------------------------
# Table 2 ("Clinical/Surgical Characteristics") - T-staging section
#   Row 34 = "T-staging" (header row, p-value cell)
#   Row 35 = "    T1"
#   Row 36 = "    T2"
#   Row 37 = "    T3"                  -> relabeled to "    T4"
#   Row 38 = "    T4"                  -> relabeled to "    Unknown/Unspecified"
#   Row 39 = "    Unknown/Unspecified" -> entire row removed
# ------------------------------------------------------------------
$t2 = $d.Tables.Item(2)

# T-staging p-value
$t2.Cell(34, 5).Range.Text = "0.006"

# T1 row
$t2.Cell(35, 2).Range.Text = "96 (91%)"
$t2.Cell(35, 3).Range.Text = "49 (94%)"
$t2.Cell(35, 4).Range.Text = "47 (87%)"

# T2 row (first value column unchanged)
$t2.Cell(36, 3).Range.Text = "0 (0%)"
$t2.Cell(36, 4).Range.Text = "6 (11%)"

# Former "T3" row becomes "T4"
$t2.Cell(37, 1).Range.Text = "    T4"
$t2.Cell(37, 2).Range.Text = "1 (0.9%)"
$t2.Cell(37, 4).Range.Text = "1 (1.9%)"

# Former "T4" row becomes "Unknown/Unspecified"
$t2.Cell(38, 1).Range.Text = "    Unknown/Unspecified"
$t2.Cell(38, 2).Range.Text = "3 (2.8%)"
$t2.Cell(38, 3).Range.Text = "3 (5.8%)"
$t2.Cell(38, 4).Range.Text = "0 (0%)"

# Remove the old trailing "Unknown/Unspecified" row entirely
$t2.Rows.Item(39).Delete()
